$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-4 values
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 10.4

$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 9.3000000000000007

$ws.Range("C4").Value = 1.3

# Row 5 ("theta_threshold_range") is removed entirely; the old row 6
# ("pie_threshold_range") shifts up to become the new row 5 with updated
# values.
$ws.Rows("5").Delete()

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 15

# Update the active cell selection to match the post-edit state.
$ws.Range("A5:XFD5").Select()

# Page setup (paper size / orientation) as in the target workbook.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
